# Add in test for objectives and endpoints
#
# Inserts a new "3.1 Primary Objectives" row into the studyDesignContent
# sheet (between the existing "3" / TRIAL OBJECTIVES... row and the
# "4" / TRIAL DESIGN row), containing a small HTML table referencing the
# Objective/Endpoint USDM attributes. Also nudges the font size used by
# the existing usdm:ref rich-text snippet in D2, and restores the two
# worksheet selections that were left behind by the editing session.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. studyDesignContent: insert the new "3.1 Primary Objectives" row
# ------------------------------------------------------------------
$content = $wb.Worksheets.Item("studyDesignContent")
$content.Activate()

# Push row 10 ("4" / TRIAL DESIGN) and everything below it down by one.
$content.Rows.Item(10).Insert()

# Section number - typed as text (looks numeric, so force text like a
# user would by leading with an apostrophe) to avoid Excel coercing it
# to a floating point value.
$content.Cells.Item(10, 1).Value = "'3.1"

# Section title.
$content.Cells.Item(10, 3).Value = "Primary Objectives"

# Section body: small HTML table wiring up the primary objective /
# endpoint USDM references.
$table = "<table>`n  <tr>`n    <th style=""vertical-align: top"">Primary Objective</th>`n    <th style=""vertical-align: top"">Primary Endpoint</th>`n  </tr>`n  <tr>`n    <td style=""vertical-align: top""><usdm:ref klass=""Objective"" id=""Objective_1"" attribute=""description""/></td>`n    <td style=""vertical-align: top""><usdm:ref klass=""Endpoint"" id=""Endpoint_1"" attribute=""description""/></td>`n  </tr>`n</table>"
$content.Cells.Item(10, 4).Value = $table

# Match the row height the content settled on in the authored workbook.
$content.Rows.Item(10).RowHeight = 113

# ------------------------------------------------------------------
# 2. Shrink the font used by the existing usdm:ref snippet in D2 from
#    12pt to 11pt (the run-level override inside the cell's text).
# ------------------------------------------------------------------
$d2 = $content.Cells.Item(2, 4)
$refRun = $d2.Characters(2, $d2.Characters(2, 200).Text.Length)
$refRun.Font.Size = 11

# ------------------------------------------------------------------
# 3. Leftover UI selections from the editing session.
# ------------------------------------------------------------------
$oe = $wb.Worksheets.Item("studyDesignOE")
$oe.Activate()
$oe.Range("C9").Select()

$content.Activate()
$content.Range("D10").Select()
